$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 ---
$ws.Range("F2").Value = 31
$ws.Range("G2").Value = "adam"
$ws.Range("I2").Value = 64
$ws.Range("J2").Value = 28.48525454334624
$ws.Range("K2").Value = 1107.651500298747
$ws.Range("L2").Value = 33.28139871307616
$ws.Range("M2").Value = 0.1577134829264187

# --- New row 3 ---
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "RNN"
$ws.Range("D3").Value = 30
$ws.Range("E3").Value = 60
$ws.Range("F3").Value = 31
$ws.Range("G3").Value = "<keras.src.optimizers.legacy.adam.Adam object at 0x79ca7e1d1600>"
$ws.Range("H3").Value = 100
$ws.Range("I3").Value = 32
$ws.Range("J3").Value = 17.52000599454657
$ws.Range("K3").Value = 463.9815377167454
$ws.Range("L3").Value = 21.54023067928349
$ws.Range("M3").Value = 0.1056961154212287

# --- New row 4 ---
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = "RNN"
$ws.Range("D4").Value = 40
$ws.Range("E4").Value = 60
$ws.Range("F4").Value = 31
$ws.Range("G4").Value = "adam"
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 64
$ws.Range("J4").Value = 32.01491631221195
$ws.Range("K4").Value = 1288.825281625744
$ws.Range("L4").Value = 35.90021283538225
$ws.Range("M4").Value = 0.1831396220764603

# --- New row 5 ---
$ws.Range("A2").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = "RNN"
$ws.Range("D5").Value = 40
$ws.Range("E5").Value = 60
$ws.Range("F5").Value = 31
$ws.Range("G5").Value = "<keras.src.optimizers.legacy.adam.Adam object at 0x79ca341ea5f0>"
$ws.Range("H5").Value = 100
$ws.Range("I5").Value = 32
$ws.Range("J5").Value = 15.59450019856665
$ws.Range("K5").Value = 421.2954551657394
$ws.Range("L5").Value = 20.52548306778039
$ws.Range("M5").Value = 0.09841049196043442
